$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44181
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 1300
$ws.Range("M2").Value = 1300
$ws.Range("P2").Value = 1300
$ws.Range("D3").Value = 44181
$ws.Range("J3").Value = 900
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 900
$ws.Range("M3").Value = 900
$ws.Range("P3").Value = 900
$ws.Range("D4").Value = 44172
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 1300
$ws.Range("L4").Value = 1300
$ws.Range("M4").Value = 1300
$ws.Range("P4").Value = 1300
$ws.Range("D5").Value = 44172
$ws.Range("J5").Value = 550
$ws.Range("D8").Value = 44165
$ws.Range("J8").Value = 720
$ws.Range("L8").Value = 1200
$ws.Range("M8").Value = 1200
$ws.Range("P8").Value = 1200
$ws.Range("D9").Value = 44165
$ws.Range("J9").Value = 750
$ws.Range("D10").Value = 44159
$ws.Range("J10").Value = 1100
$ws.Range("D11").Value = 44159
$ws.Range("J11").Value = 800
$ws.Range("D12").Value = 44169
$ws.Range("J12").Value = 950
$ws.Range("D13").Value = 44169
$ws.Range("J13").Value = 800
$ws.Range("K13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = 1000
$ws.Range("P13").Value = 1000
$ws.Range("D14").Value = 44179
$ws.Range("J14").Value = 980
$ws.Range("K14").Value = 1200
$ws.Range("L14").Value = 1200
$ws.Range("M14").Value = 1200
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 1200
$ws.Range("D15").Value = 44174
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 2800
$ws.Range("K15").Value = 1200
$ws.Range("L15").Value = 1250
$ws.Range("M15").Value = 1221
$ws.Range("P15").Value = 1221
$ws.Range("D16").Value = 44174
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 1300
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 1000
$ws.Range("P16").Value = 1000
$ws.Range("D17").Value = 44175
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 1500
$ws.Range("K17").Value = 1300
$ws.Range("L17").Value = 1300
$ws.Range("M17").Value = 1300
$ws.Range("P17").Value = 1300
$ws.Range("D18").Value = 44175
$ws.Range("I18").Value = "Segunda"
$ws.Range("J18").Value = 1450
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = 1000
$ws.Range("P18").Value = 1000
$ws.Range("D19").Value = 44161
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 1600
$ws.Range("K19").Value = 1300
$ws.Range("L19").Value = 1300
$ws.Range("M19").Value = 1300
$ws.Range("P19").Value = 1300
$ws.Range("D20").Value = 44161
$ws.Range("I20").Value = "Segunda"
$ws.Range("J20").Value = 1850
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = 1000
$ws.Range("P20").Value = 1000
$ws.Range("D21").Value = 44167
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 1430
$ws.Range("K21").Value = 1200
$ws.Range("L21").Value = 1300
$ws.Range("M21").Value = 1248
$ws.Range("P21").Value = 1248
$ws.Range("D22").Value = 44167
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 1000
$ws.Range("P22").Value = 1000
$ws.Range("D23").Value = 44162
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 1200
$ws.Range("K23").Value = 1300
$ws.Range("L23").Value = 1300
$ws.Range("M23").Value = 1300
$ws.Range("P23").Value = 1300
$ws.Range("I24").Value = "Segunda"
$ws.Range("J24").Value = 800
$ws.Range("K24").Value = 1000
$ws.Range("L24").Value = 1000
$ws.Range("M24").Value = 1000
$ws.Range("P24").Value = 1000
$ws.Range("D25").Value = 44176
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 2500
$ws.Range("K25").Value = 1200
$ws.Range("L25").Value = 1300
$ws.Range("M25").Value = 1256
$ws.Range("P25").Value = 1256
$ws.Range("D26").Value = 44176
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 1500
$ws.Range("K26").Value = 1000
$ws.Range("L26").Value = 1000
$ws.Range("M26").Value = 1000
$ws.Range("P26").Value = 1000
$ws.Range("D27").Value = 44160
$ws.Range("I27").Value = "Primera"
$ws.Range("K27").Value = 1300
$ws.Range("L27").Value = 1300
$ws.Range("M27").Value = 1300
$ws.Range("P27").Value = 1300
$ws.Range("D28").Value = 44160
$ws.Range("I28").Value = "Segunda"
$ws.Range("J28").Value = 850
$ws.Range("K28").Value = 1000
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 1000
$ws.Range("O28").Value = "Provincia de Quillota"
$ws.Range("P28").Value = 1000
